$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mislabeled header: FT_Goals_H -> FT_Goals_A (column G)
$ws.Range("G1").Value = "FT_Goals_A"

# Append new match rows (178-190)
# Row 178
$ws.Range("A178").Value = 'E0'
$ws.Range("B178").Value = '''12/01/2023'
$ws.Range("C178").Value = '20:00'
$ws.Range("D178").Value = 'Fulham'
$ws.Range("E178").Value = 'Chelsea'
$ws.Range("F178").Value = 2
$ws.Range("G178").Value = 1
$ws.Range("H178").Value = 'H'
$ws.Range("I178").Value = 1
$ws.Range("J178").Value = 0
$ws.Range("K178").Value = 'H'
$ws.Range("L178").Value = 'D Coote'
$ws.Range("M178").Value = 8
$ws.Range("N178").Value = 20
$ws.Range("O178").Value = 3
$ws.Range("P178").Value = 10
$ws.Range("Q178").Value = 12
$ws.Range("R178").Value = 16
$ws.Range("S178").Value = 5
$ws.Range("T178").Value = 7
$ws.Range("U178").Value = 4
$ws.Range("V178").Value = 3
$ws.Range("W178").Value = 0
$ws.Range("X178").Value = 1
$ws.Range("Y178").Value = 3.3
$ws.Range("Z178").Value = 3.5
$ws.Range("AA178").Value = 2.15
$ws.Range("AB178").Value = 1.87
$ws.Range("AC178").Value = 2.03

# Row 179
$ws.Range("A179").Value = 'E0'
$ws.Range("B179").Value = '13/01/2023'
$ws.Range("C179").Value = '20:00'
$ws.Range("D179").Value = 'Aston Villa'
$ws.Range("E179").Value = 'Leeds'
$ws.Range("F179").Value = 2
$ws.Range("G179").Value = 1
$ws.Range("H179").Value = 'H'
$ws.Range("I179").Value = 1
$ws.Range("J179").Value = 0
$ws.Range("K179").Value = 'H'
$ws.Range("L179").Value = 'M Oliver'
$ws.Range("M179").Value = 11
$ws.Range("N179").Value = 16
$ws.Range("O179").Value = 5
$ws.Range("P179").Value = 4
$ws.Range("Q179").Value = 12
$ws.Range("R179").Value = 8
$ws.Range("S179").Value = 0
$ws.Range("T179").Value = 11
$ws.Range("U179").Value = 2
$ws.Range("V179").Value = 2
$ws.Range("W179").Value = 0
$ws.Range("X179").Value = 0
$ws.Range("Y179").Value = 1.91
$ws.Range("Z179").Value = 3.6
$ws.Range("AA179").Value = 4
$ws.Range("AB179").Value = 1.8
$ws.Range("AC179").Value = 2

# Row 180
$ws.Range("A180").Value = 'E0'
$ws.Range("B180").Value = '14/01/2023'
$ws.Range("C180").Value = '12:30'
$ws.Range("D180").Value = 'Man United'
$ws.Range("E180").Value = 'Man City'
$ws.Range("F180").Value = 2
$ws.Range("G180").Value = 1
$ws.Range("H180").Value = 'H'
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = 0
$ws.Range("K180").Value = 'D'
$ws.Range("L180").Value = 'S Attwell'
$ws.Range("M180").Value = 8
$ws.Range("N180").Value = 5
$ws.Range("O180").Value = 4
$ws.Range("P180").Value = 1
$ws.Range("Q180").Value = 15
$ws.Range("R180").Value = 5
$ws.Range("S180").Value = 1
$ws.Range("T180").Value = 2
$ws.Range("U180").Value = 3
$ws.Range("V180").Value = 0
$ws.Range("W180").Value = 0
$ws.Range("X180").Value = 0
$ws.Range("Y180").Value = 4.2
$ws.Range("Z180").Value = 3.6
$ws.Range("AA180").Value = 1.85
$ws.Range("AB180").Value = 1.73
$ws.Range("AC180").Value = 2.1

# Row 181
$ws.Range("A181").Value = 'E0'
$ws.Range("B181").Value = '14/01/2023'
$ws.Range("C181").Value = '15:00'
$ws.Range("D181").Value = 'Brighton'
$ws.Range("E181").Value = 'Liverpool'
$ws.Range("F181").Value = 3
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 'H'
$ws.Range("I181").Value = 0
$ws.Range("J181").Value = 0
$ws.Range("K181").Value = 'D'
$ws.Range("L181").Value = 'D England'
$ws.Range("M181").Value = 16
$ws.Range("N181").Value = 6
$ws.Range("O181").Value = 9
$ws.Range("P181").Value = 2
$ws.Range("Q181").Value = 8
$ws.Range("R181").Value = 15
$ws.Range("S181").Value = 7
$ws.Range("T181").Value = 1
$ws.Range("U181").Value = 1
$ws.Range("V181").Value = 3
$ws.Range("W181").Value = 0
$ws.Range("X181").Value = 0
$ws.Range("Y181").Value = 3.4
$ws.Range("Z181").Value = 3.6
$ws.Range("AA181").Value = 2.1
$ws.Range("AB181").Value = 1.57
$ws.Range("AC181").Value = 2.38

# Row 182
$ws.Range("A182").Value = 'E0'
$ws.Range("B182").Value = '14/01/2023'
$ws.Range("C182").Value = '15:00'
$ws.Range("D182").Value = 'Everton'
$ws.Range("E182").Value = 'Southampton'
$ws.Range("F182").Value = 1
$ws.Range("G182").Value = 2
$ws.Range("H182").Value = 'A'
$ws.Range("I182").Value = 1
$ws.Range("J182").Value = 0
$ws.Range("K182").Value = 'H'
$ws.Range("L182").Value = 'J Brooks'
$ws.Range("M182").Value = 12
$ws.Range("N182").Value = 13
$ws.Range("O182").Value = 4
$ws.Range("P182").Value = 5
$ws.Range("Q182").Value = 7
$ws.Range("R182").Value = 19
$ws.Range("S182").Value = 6
$ws.Range("T182").Value = 4
$ws.Range("U182").Value = 0
$ws.Range("V182").Value = 3
$ws.Range("W182").Value = 0
$ws.Range("X182").Value = 0
$ws.Range("Y182").Value = 2.25
$ws.Range("Z182").Value = 3.25
$ws.Range("AA182").Value = 3.3
$ws.Range("AB182").Value = 2.1
$ws.Range("AC182").Value = 1.73

# Row 183
$ws.Range("A183").Value = 'E0'
$ws.Range("B183").Value = '14/01/2023'
$ws.Range("C183").Value = '15:00'
$ws.Range("D183").Value = 'Nott''m Forest'
$ws.Range("E183").Value = 'Leicester'
$ws.Range("F183").Value = 2
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 'H'
$ws.Range("I183").Value = 0
$ws.Range("J183").Value = 0
$ws.Range("K183").Value = 'D'
$ws.Range("L183").Value = 'P Tierney'
$ws.Range("M183").Value = 15
$ws.Range("N183").Value = 6
$ws.Range("O183").Value = 5
$ws.Range("P183").Value = 1
$ws.Range("Q183").Value = 11
$ws.Range("R183").Value = 12
$ws.Range("S183").Value = 3
$ws.Range("T183").Value = 5
$ws.Range("U183").Value = 1
$ws.Range("V183").Value = 2
$ws.Range("W183").Value = 0
$ws.Range("X183").Value = 0
$ws.Range("Y183").Value = 2.7
$ws.Range("Z183").Value = 3.3
$ws.Range("AA183").Value = 2.63
$ws.Range("AB183").Value = 2.03
$ws.Range("AC183").Value = 1.87

# Row 184
$ws.Range("A184").Value = 'E0'
$ws.Range("B184").Value = '14/01/2023'
$ws.Range("C184").Value = '15:00'
$ws.Range("D184").Value = 'Wolves'
$ws.Range("E184").Value = 'West Ham'
$ws.Range("F184").Value = 1
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 'H'
$ws.Range("I184").Value = 0
$ws.Range("J184").Value = 0
$ws.Range("K184").Value = 'D'
$ws.Range("L184").Value = 'S Hooper'
$ws.Range("M184").Value = 17
$ws.Range("N184").Value = 16
$ws.Range("O184").Value = 4
$ws.Range("P184").Value = 4
$ws.Range("Q184").Value = 10
$ws.Range("R184").Value = 7
$ws.Range("S184").Value = 4
$ws.Range("T184").Value = 3
$ws.Range("U184").Value = 0
$ws.Range("V184").Value = 2
$ws.Range("W184").Value = 0
$ws.Range("X184").Value = 0
$ws.Range("Y184").Value = 2.63
$ws.Range("Z184").Value = 3.3
$ws.Range("AA184").Value = 2.7
$ws.Range("AB184").Value = 2.2
$ws.Range("AC184").Value = 1.67

# Row 185
$ws.Range("A185").Value = 'E0'
$ws.Range("B185").Value = '14/01/2023'
$ws.Range("C185").Value = '17:30'
$ws.Range("D185").Value = 'Brentford'
$ws.Range("E185").Value = 'Bournemouth'
$ws.Range("F185").Value = 2
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 'H'
$ws.Range("I185").Value = 1
$ws.Range("J185").Value = 0
$ws.Range("K185").Value = 'H'
$ws.Range("L185").Value = 'J Gillett'
$ws.Range("M185").Value = 12
$ws.Range("N185").Value = 7
$ws.Range("O185").Value = 4
$ws.Range("P185").Value = 2
$ws.Range("Q185").Value = 5
$ws.Range("R185").Value = 9
$ws.Range("S185").Value = 5
$ws.Range("T185").Value = 3
$ws.Range("U185").Value = 1
$ws.Range("V185").Value = 1
$ws.Range("W185").Value = 0
$ws.Range("X185").Value = 0
$ws.Range("Y185").Value = 1.65
$ws.Range("Z185").Value = 4
$ws.Range("AA185").Value = 5.25
$ws.Range("AB185").Value = 1.8
$ws.Range("AC185").Value = 2

# Row 186
$ws.Range("A186").Value = 'E0'
$ws.Range("B186").Value = '15/01/2023'
$ws.Range("C186").Value = '14:00'
$ws.Range("D186").Value = 'Chelsea'
$ws.Range("E186").Value = 'Crystal Palace'
$ws.Range("F186").Value = 1
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 'H'
$ws.Range("I186").Value = 0
$ws.Range("J186").Value = 0
$ws.Range("K186").Value = 'D'
$ws.Range("L186").Value = 'P Bankes'
$ws.Range("M186").Value = 15
$ws.Range("N186").Value = 10
$ws.Range("O186").Value = 5
$ws.Range("P186").Value = 5
$ws.Range("Q186").Value = 10
$ws.Range("R186").Value = 17
$ws.Range("S186").Value = 11
$ws.Range("T186").Value = 7
$ws.Range("U186").Value = 2
$ws.Range("V186").Value = 5
$ws.Range("W186").Value = 0
$ws.Range("X186").Value = 0
$ws.Range("Y186").Value = 1.62
$ws.Range("Z186").Value = 4
$ws.Range("AA186").Value = 5.25
$ws.Range("AB186").Value = 2.02
$ws.Range("AC186").Value = 1.88

# Row 187
$ws.Range("A187").Value = 'E0'
$ws.Range("B187").Value = '15/01/2023'
$ws.Range("C187").Value = '14:00'
$ws.Range("D187").Value = 'Newcastle'
$ws.Range("E187").Value = 'Fulham'
$ws.Range("F187").Value = 1
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 'H'
$ws.Range("I187").Value = 0
$ws.Range("J187").Value = 0
$ws.Range("K187").Value = 'D'
$ws.Range("L187").Value = 'R Jones'
$ws.Range("M187").Value = 20
$ws.Range("N187").Value = 5
$ws.Range("O187").Value = 5
$ws.Range("P187").Value = 0
$ws.Range("Q187").Value = 6
$ws.Range("R187").Value = 14
$ws.Range("S187").Value = 10
$ws.Range("T187").Value = 5
$ws.Range("U187").Value = 1
$ws.Range("V187").Value = 3
$ws.Range("W187").Value = 0
$ws.Range("X187").Value = 0
$ws.Range("Y187").Value = 1.5
$ws.Range("Z187").Value = 4.33
$ws.Range("AA187").Value = 6.5
$ws.Range("AB187").Value = 1.73
$ws.Range("AC187").Value = 2.1

# Row 188
$ws.Range("A188").Value = 'E0'
$ws.Range("B188").Value = '15/01/2023'
$ws.Range("C188").Value = '16:30'
$ws.Range("D188").Value = 'Tottenham'
$ws.Range("E188").Value = 'Arsenal'
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 2
$ws.Range("H188").Value = 'A'
$ws.Range("I188").Value = 0
$ws.Range("J188").Value = 2
$ws.Range("K188").Value = 'A'
$ws.Range("L188").Value = 'C Pawson'
$ws.Range("M188").Value = 17
$ws.Range("N188").Value = 14
$ws.Range("O188").Value = 7
$ws.Range("P188").Value = 5
$ws.Range("Q188").Value = 16
$ws.Range("R188").Value = 15
$ws.Range("S188").Value = 4
$ws.Range("T188").Value = 3
$ws.Range("U188").Value = 4
$ws.Range("V188").Value = 2
$ws.Range("W188").Value = 0
$ws.Range("X188").Value = 0
$ws.Range("Y188").Value = 3.1
$ws.Range("Z188").Value = 3.6
$ws.Range("AA188").Value = 2.2
$ws.Range("AB188").Value = 1.8
$ws.Range("AC188").Value = 2

# Row 189
$ws.Range("A189").Value = 'E0'
$ws.Range("B189").Value = '18/01/2023'
$ws.Range("C189").Value = '20:00'
$ws.Range("D189").Value = 'Crystal Palace'
$ws.Range("E189").Value = 'Man United'
$ws.Range("F189").Value = 1
$ws.Range("G189").Value = 1
$ws.Range("H189").Value = 'D'
$ws.Range("I189").Value = 0
$ws.Range("J189").Value = 1
$ws.Range("K189").Value = 'A'
$ws.Range("L189").Value = 'R Jones'
$ws.Range("M189").Value = 10
$ws.Range("N189").Value = 15
$ws.Range("O189").Value = 5
$ws.Range("P189").Value = 4
$ws.Range("Q189").Value = 9
$ws.Range("R189").Value = 10
$ws.Range("S189").Value = 3
$ws.Range("T189").Value = 3
$ws.Range("U189").Value = 1
$ws.Range("V189").Value = 2
$ws.Range("W189").Value = 0
$ws.Range("X189").Value = 0
$ws.Range("Y189").Value = 4.5
$ws.Range("Z189").Value = 3.6
$ws.Range("AA189").Value = 1.8
$ws.Range("AB189").Value = 1.97
$ws.Range("AC189").Value = 1.93

# Row 190
$ws.Range("A190").Value = 'E0'
$ws.Range("B190").Value = '19/01/2023'
$ws.Range("C190").Value = '20:00'
$ws.Range("D190").Value = 'Man City'
$ws.Range("E190").Value = 'Tottenham'
$ws.Range("F190").Value = 4
$ws.Range("G190").Value = 2
$ws.Range("H190").Value = 'H'
$ws.Range("I190").Value = 0
$ws.Range("J190").Value = 2
$ws.Range("K190").Value = 'A'
$ws.Range("L190").Value = 'S Hooper'
$ws.Range("M190").Value = 16
$ws.Range("N190").Value = 9
$ws.Range("O190").Value = 6
$ws.Range("P190").Value = 3
$ws.Range("Q190").Value = 10
$ws.Range("R190").Value = 14
$ws.Range("S190").Value = 8
$ws.Range("T190").Value = 3
$ws.Range("U190").Value = 1
$ws.Range("V190").Value = 2
$ws.Range("W190").Value = 0
$ws.Range("X190").Value = 0
$ws.Range("Y190").Value = 1.36
$ws.Range("Z190").Value = 5
$ws.Range("AA190").Value = 8.5
$ws.Range("AB190").Value = 1.53
$ws.Range("AC190").Value = 2.5

